$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 - section header (bold, like A1/A10/A15/A20/A25)
$ws.Range("A32").Value = "H2O2"
$ws.Range("A32").Font.Bold = $true

# Row 33
$ws.Range("A33").Value = "init conc"
$ws.Range("B33").Value = "w."
$ws.Range("C33").Value = 0.3

# Row 34
$ws.Range("A34").Value = "density"
$ws.Range("B34").Value = "g/ml"
$ws.Range("C34").Value = 1.11

# Row 35
$ws.Range("A35").Value = "molar w."
$ws.Range("B35").Value = "g/mol"
$ws.Range("C35").Value = 34.01468

# Row 36
$ws.Range("A36").Value = "init conc"
$ws.Range("B36").Value = "M"
$ws.Range("C36").Formula = "=C33*C34/C35*1000"

# Row 37 left blank (but must materialize as an empty row, matching rows 9/14/19/24)
$ws.Rows.Item(37).RowHeight = 14.25

# Row 38
$ws.Range("A38").Value = "conc needed"
$ws.Range("B38").Value = "μM"
$ws.Range("C38").Value = 10000

# Row 39
$ws.Range("A39").Value = "V"
$ws.Range("B39").Value = "ml"
$ws.Range("C39").Formula = "=C38*10^(-6)/C36*1000"

# Row 40 left blank (but must materialize as an empty row)
$ws.Rows.Item(40).RowHeight = 14.25

# New column-G width entry (matching the diff's ~10.85-char bestFit width;
# the host quantizes to whole pixels, so 10.0 is the closest achievable input)
$ws.Columns.Item(7).ColumnWidth = 10.0
